$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings (e.g. "1.003", "6.710",
# "25.682.94"). Pre-set NumberFormat to Text ("@") so the COM layer
# stores the literal digits instead of auto-converting to a number
# (which would drop trailing zeros / mis-parse the dotted thousands).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.682.94'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '1.620.59'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '214.59'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '0.5076'
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.2559'
$ws.Range("E8").Value = '  -1.76%  '
$ws.Range("D9").Value = '0.06351'
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("D10").Value = '19.19'
$ws.Range("E10").Value = '  -4.03%  '
$ws.Range("D11").Value = '0.07759'
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = '4.227'
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("D13").Value = '1.626.33'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").Value = '1.844.60'
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").Value = '0.5535'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").Value = '63.41'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").Value = '0.0₅7497'
$ws.Range("E17").Value = '  -3.88%  '
$ws.Range("D18").Value = '25.715.37'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '193.19'
$ws.Range("E20").Value = '  -3.67%  '
$ws.Range("D21").Value = '4.358'
$ws.Range("E21").Value = '  -2.99%  '
$ws.Range("D22").Value = '9.733'
$ws.Range("E22").Value = '  -3.15%  '
$ws.Range("D23").Value = '5.948'
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '1.849'
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").Value = '140.39'
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("D27").Value = '0.1251'
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("D28").Value = '6.710'
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("D29").Value = '15.41'
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("D31").Value = '0.04848'
$ws.Range("E31").Value = '  -2.12%  '
$ws.Range("D32").Value = '3.282'
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("D33").Value = '3.162'
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("D34").Value = '1.538'
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").Value = '2.366'
$ws.Range("D36").Value = '0.8901'
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("D37").Value = '1.124.20'
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").Value = '2.530'
$ws.Range("E38").Value = '  -2.36%  '
$ws.Range("D39").Value = '0.5468'
$ws.Range("E39").Value = '  -2.79%  '
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").Value = '5.566'
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("D43").Value = '0.7926'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = '96.89'
$ws.Range("E44").Value = '  -3.16%  '
$ws.Range("D45").Value = '1.769.12'
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("D46").Value = '0.0₈112'
$ws.Range("E46").Value = '  -9.59%  '
$ws.Range("D47").Value = '0.4411'
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").Value = '0.05109'
$ws.Range("E49").Value = '  -3.26%  '
$ws.Range("D50").Value = '7.572'
$ws.Range("E50").Value = '  +1.42%  '
$ws.Range("D51").Value = '0.9973'
$ws.Range("E51").Value = '  -0.85%  '

# Reset to the default Normal style so the saved cells do not carry an
# explicit visible number-format, matching the original workbook where
# these cells had no style override.
$ws.Range("D2:D51").Style = "Normal"
